# "working first DEMO : 7inch Touch Screen + STM32F429ZI"
#
# Typography sheet: remove the now-unused "Typography_00" typography
# definition row (row 7, columns B:J) - its associated font
# (MATURASC.TTF) is no longer referenced.
#
# Translation sheet: replace the single "HELLO WORLD !" translation
# (which previously spanned two lines and was left/center aligned with a
# throw-away id) with two simple left-aligned "HELLO WORLD !" entries
# (SingleUseId1 and SingleUseId2) ready for the 7" touch screen demo.

$wb = $excel.ActiveWorkbook

$wsTypography = $wb.Worksheets.Item("Typography")
$wsTranslation = $wb.Worksheets.Item("Translation")

# --- Typography sheet: drop row 7 ("Typography_00") ---
$wsTypography.Range("B7:J7").ClearContents()

# --- Translation sheet: rewrite rows 4 and 5 ---
$wsTranslation.Range("B4").Value = "SingleUseId1"
$wsTranslation.Range("C4").Value = "Large"
$wsTranslation.Range("D4").Value = "Left"
$wsTranslation.Range("E4").Value = "LTR"
$wsTranslation.Range("F4").Value = "HELLO WORLD !"

$wsTranslation.Range("B5").Value = "SingleUseId2"
$wsTranslation.Range("C5").Value = "Large"
$wsTranslation.Range("D5").Value = "Left"
$wsTranslation.Range("E5").Value = "LTR"
$wsTranslation.Range("F5").Value = "HELLO WORLD !"

# Rows 6 and 7 previously existed only as empty placeholder rows; they are
# already empty, so nothing further needs to be cleared there.

Write-Host "Applied 7-inch touch screen demo text edits."
